$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D2").Value = "2016-01-26 05:39:51"
$wsZh.Range("G2").Value = "2016-01-26 05:40:34"

$wsDe.Range("D2").Value = "2016-01-26 05:40:01"
$wsDe.Range("G2").Value = "2016-01-26 05:40:51"
